$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header renames
$ws.Range("E1").Value = "risk"
$ws.Range("F1").Value = "ln_tw_risk"

# Data value updates (columns E and F)
$ws.Range("E2").Value = 0.00614137489401504
$ws.Range("F2").Value = 5.22647841858601
$ws.Range("E4").Value = 0.0127563015870867
$ws.Range("E5").Value = 0.209665059032301
$ws.Range("F5").Value = 8.999907763926535
$ws.Range("E6").Value = 0.0227748818358775
$ws.Range("E7").Value = 0.168347664564615
$ws.Range("F7").Value = 4.882978316844667
$ws.Range("E8").Value = 0.0256975843456798
$ws.Range("F8").Value = 2.797961673440694
$ws.Range("E9").Value = 0.00272292284457534
$ws.Range("F9").Value = 3.18077305275393
$ws.Range("E10").Value = 0.0364116937993349
$ws.Range("F10").Value = 2.232989030073072
$ws.Range("E11").Value = 0.288755608721611
$ws.Range("F11").Value = 7.284534584519009
$ws.Range("E12").Value = 0.0186017066066913
$ws.Range("F12").Value = 4.025763996142878
$ws.Range("E13").Value = 0.214874128821924
$ws.Range("F13").Value = 6.251424486977098
$ws.Range("E14").Value = 1.7381126844623
$ws.Range("F14").Value = 10.39131042825277
$ws.Range("E15").Value = 0.193823271504487
$ws.Range("F15").Value = 5.615240299123142
$ws.Range("E16").Value = 0.0583578599868634
$ws.Range("F16").Value = 4.007962236780179
$ws.Range("E17").Value = 0.425575082414661
$ws.Range("F17").Value = 8.16371069703683
$ws.Range("E18").Value = 0.782980903505756
$ws.Range("F18").Value = 8.667462307344751
$ws.Range("E19").Value = 0.00670371765539473
$ws.Range("F19").Value = 2.571630705120837
$ws.Range("E20").Value = 0.13908364323229
$ws.Range("F20").Value = 6.459945022965127
$ws.Range("E21").Value = 1.02003057516788
$ws.Range("F21").Value = 10.00183013011388
$ws.Range("E22").Value = [double]"2.21977405807772e-05"
$ws.Range("E23").Value = 1.4112805513978
$ws.Range("F23").Value = 10.23269209319368
$ws.Range("E24").Value = 0.029530394219294
$ws.Range("E25").Value = [double]"8.8790962323109e-05"
$ws.Range("E26").Value = 0.145025238461078
$ws.Range("F26").Value = 5.510705086771525
$ws.Range("E28").Value = 0.0166039099544214
$ws.Range("F28").Value = 4.306927013280672
$ws.Range("E29").Value = 0.0949767326982856
$ws.Range("F29").Value = 5.797843379114219
$ws.Range("E30").Value = 0.0513951686913596
$ws.Range("F30").Value = 5.033724635981492
$ws.Range("E31").Value = 0.0394379857651809
$ws.Range("F31").Value = 5.691169349051433
$ws.Range("E32").Value = 0.211714650412593
$ws.Range("F32").Value = 7.623928339781987
$ws.Range("E33").Value = 1.13150762836454
$ws.Range("F33").Value = 9.599867261050562
$ws.Range("E34").Value = 0.0473181836713568
$ws.Range("F34").Value = 4.523569600552934
$ws.Range("E37").Value = 0.00937484577194826
$ws.Range("F37").Value = 3.646917334359665
$ws.Range("E38").Value = 0.0145839155615707
$ws.Range("F38").Value = 3.272575131061424
$ws.Range("E39").Value = 0.00785800016559514
$ws.Range("F39").Value = 1.799586756102128
$ws.Range("E40").Value = 0.125594816206038
$ws.Range("F40").Value = 6.696172979167821
$ws.Range("E42").Value = 0.00318167614991141
$ws.Range("E43").Value = 0.00685170259259991
$ws.Range("F43").Value = 1.79858780557809
$ws.Range("E44").Value = 0.0108990906251616
$ws.Range("F44").Value = 3.59435906206405
$ws.Range("E45").Value = 0.00573441631670079
$ws.Range("E46").Value = 0.331937613398083
$ws.Range("F46").Value = 9.404969460413254
$ws.Range("E47").Value = 0.0552723740461353
$ws.Range("F47").Value = 6.536906259515253
$ws.Range("E48").Value = 0.127052467837509
$ws.Range("F48").Value = 6.199538984382349
$ws.Range("E49").Value = 0.138388114027426
$ws.Range("F49").Value = 7.098463705135256
$ws.Range("E51").Value = 0.0410732193212982
$ws.Range("F51").Value = 3.703813768852255
$ws.Range("E52").Value = 1.11930627029197
$ws.Range("F52").Value = 6.888815858224056
$ws.Range("E54").Value = 0.0750949563847694
$ws.Range("E55").Value = 0.0008583126357900529
$ws.Range("F55").Value = 1.610295856930237
$ws.Range("E56").Value = 0.00745844083514115
$ws.Range("F56").Value = 4.014763949428953
$ws.Range("E57").Value = 0.148828451347251
$ws.Range("F57").Value = 7.264829958449162
$ws.Range("E58").Value = 0.000118387949764145
$ws.Range("E59").Value = 0.0215466068570744
$ws.Range("F59").Value = 4.888852210590512
$ws.Range("E60").Value = 0.491975923738626
$ws.Range("F60").Value = 9.617218721002168
$ws.Range("E61").Value = 0.127821989510976
$ws.Range("F61").Value = 2.685237687366762
$ws.Range("E63").Value = 0
$ws.Range("E64").Value = [double]"8.8790962323109e-05"
$ws.Range("E65").Value = 0
